$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.176578283309937
$ws.Range("B1").Value = 4.545527935028076
$ws.Range("C1").Value = 3.96592903137207
$ws.Range("D1").Value = 1.509686470031738
$ws.Range("E1").Value = 0.8882495164871216
